# Fill in row 20 of the "Cliente" sheet with the newest client/sale record
# (ANTONIO VIEIRA), which was previously an empty placeholder row.

$wb = $excel.ActiveWorkbook
$ws = $wb.Sheets.Item("Cliente")

$ws.Range("B20").Value = "ANTONIO VIEIRA"
$ws.Range("C20").Value = "533ae974ff5ab9eaed4a9f8074909ec5"
$ws.Range("D20").Value = Get-Date -Year 2022 -Month 10 -Day 18
$ws.Range("E20").Value = 365
$ws.Range("F20").Value = "-"
$ws.Range("G20").Value = "VENDA 16 (18/10)"
